# Update "want-to-go" counts (column F) and, for two sold-out rows, the
# "minimum price" column (G) text, across the four worksheets.
#
# Sheet order (tab order) in this workbook:
#   1 = 展览 (Exhibitions)
#   2 = 演出 (Shows)
#   3 = 本地生活 (Local life)
#   4 = 全部类型 (All types, aggregates the other three sheets)

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value  = 1396
$ws1.Range("F4").Value  = 13369
$ws1.Range("F5").Value  = 769
$ws1.Range("F10").Value = 1916
$ws1.Range("F13").Value = 20869
$ws1.Range("G13").Value = "已售罄"
$ws1.Range("F16").Value = 492
$ws1.Range("F26").Value = 18
$ws1.Range("F28").Value = 62
$ws1.Range("F29").Value = 379

# --- Sheet 2: 演出 -----------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F4").Value  = 4478
$ws2.Range("F11").Value = 388

# --- Sheet 3: 本地生活 -------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 907
$ws3.Range("F3").Value = 4433
$ws3.Range("F4").Value = 104

# --- Sheet 4: 全部类型 -------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 907
$ws4.Range("F5").Value  = 1396
$ws4.Range("F6").Value  = 13369
$ws4.Range("F8").Value  = 769
$ws4.Range("F9").Value  = 4433
$ws4.Range("F13").Value = 1916
$ws4.Range("F16").Value = 104
$ws4.Range("F17").Value = 20869
$ws4.Range("G17").Value = "已售罄"
$ws4.Range("F19").Value = 4478
$ws4.Range("F23").Value = 492
$ws4.Range("F29").Value = 388
$ws4.Range("F41").Value = 18
$ws4.Range("F43").Value = 62
$ws4.Range("F45").Value = 379
